# weather_parameters.xlsx - Alycia's latest edits (switch the model run's
# weather inputs from Namibia to Kenya), keeping the same code/layout so
# the run can be split consistently with Claire.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Filename parameter: Namibia-2022 -> Kenya-2022
$ws.Cells.Item(8, 2).Value = "Kenya-2022"

# Bounding box for the new (Kenya) region
$ws.Cells.Item(4, 2).Value = 34    # Minimum longitude (deg)
$ws.Cells.Item(5, 2).Value = 42    # Maximum longitude (deg)
$ws.Cells.Item(6, 2).Value = -5    # Minimum latitude (deg)
$ws.Cells.Item(7, 2).Value = 5     # Maximum latitude (deg)

# Leave the selection where the author last left it before saving
$ws.Range("F10").Select()
